$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 19575605.8673771
$ws.Range("E2").Value = 13.86384647080068
$ws.Range("G2").Value = 19575624.63053851

$ws.Range("B3").Value = 0.04172184405617529
$ws.Range("C3").Value = 0.002658071450198252
$ws.Range("D3").Value = 18.71679738969934
$ws.Range("E3").Value = 13.86384647080068
$ws.Range("G3").Value = 32.62502377600639

$ws.Range("B4").Value = 1.445647641019636
$ws.Range("C4").Value = 2919.202174992006
$ws.Range("D4").Value = 3993.344853322108
$ws.Range("E4").Value = 2797.565817734744
$ws.Range("G4").Value = 9711.558493689876
